# Character sheet: rename the "height"/"radius" float columns (C/D) to the
# new character stats "moveSpeed"/"weight", update the sample data row, and
# drop the now-unused "radius" column E entirely (the grid shrinks from
# A1:E5 to A1:D5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Character")
$ws.Activate()

# Header row (fieldName)
$ws.Range("C1").Value = "moveSpeed"
$ws.Range("D1").Value = "weight"

# Sample data row
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 50

# Remove the old "radius" column (E); C/D keep their "float" type row (row 4)
# and styles, everything to the right just shifts left.
$ws.Columns.Item(5).Delete()

$ws.Range("D4").Select()
